$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Report generated for the handback of c83ea4d2-a7e8-450c-846b-80a60074867a
# (row 7 on both the "zh-cn" and "de-de" sheets): the handback xliff has been
# received, so the "Latest Target File" column becomes a hyperlink to the
# source markdown, "Latest Handback File" / "Latest Handback DateTime" are
# populated, and an "Error Detail" note is recorded because the handback is
# based on an older version of the source file.
# ---------------------------------------------------------------------------

$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add(
    $zhcn.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3d98506a59ecf6059580dfb5de8993517f0329ff/e2e/c83ea4d2-a7e8-450c-846b-80a60074867a.md",
    "",
    "",
    "c83ea4d2-a7e8-450c-846b-80a60074867a.md"
)
$zhcn.Range("I7").Font.Underline = $true
$zhcn.Range("I7").Font.Color = 15570276

$zhcn.Range("J7").Value = "c83ea4d2-a7e8-450c-846b-80a60074867a.b07396c9bd66bc9cfbde25f7174520906c666fd9.zh-cn.xlf"
$zhcn.Range("K7").Value = "2016-08-15 10:53:09"
$zhcn.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e0e0a12686bdc4ae7e1f16a962c97706cf26872/e2e/c83ea4d2-a7e8-450c-846b-80a60074867a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d98506a59ecf6059580dfb5de8993517f0329ff/e2e/c83ea4d2-a7e8-450c-846b-80a60074867a.md."

$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add(
    $dede.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3d98506a59ecf6059580dfb5de8993517f0329ff/e2e/c83ea4d2-a7e8-450c-846b-80a60074867a.md",
    "",
    "",
    "c83ea4d2-a7e8-450c-846b-80a60074867a.md"
)
$dede.Range("I7").Font.Underline = $true
$dede.Range("I7").Font.Color = 15570276

$dede.Range("J7").Value = "c83ea4d2-a7e8-450c-846b-80a60074867a.b07396c9bd66bc9cfbde25f7174520906c666fd9.de-de.xlf"
$dede.Range("K7").Value = "2016-08-15 10:53:16"
$dede.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e0e0a12686bdc4ae7e1f16a962c97706cf26872/e2e/c83ea4d2-a7e8-450c-846b-80a60074867a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d98506a59ecf6059580dfb5de8993517f0329ff/e2e/c83ea4d2-a7e8-450c-846b-80a60074867a.md."
